function Get-LocalName {
    param($lnKey)
    $lnIdx = $lnKey.IndexOf(":")
    if ($lnIdx -ge 0) {
        return $lnKey.Substring($lnIdx + 1)
    }
    return $lnKey
}

function Reorder-Tag {
    param($rtMatch)

    $rtTagName = $rtMatch.Groups[1].Value
    $rtAttrStr = $rtMatch.Groups[2].Value
    $rtSelfClose = $rtMatch.Groups[3].Value

    if ($rtAttrStr.Trim().Length -eq 0) {
        return $rtMatch.Value
    }

    $rtAttrMatches = [regex]::Matches($rtAttrStr, '([\w:.\-]+)="((?:[^"])*)"')
    $rtNsAttrs = New-Object System.Collections.ArrayList
    $rtNormAttrs = New-Object System.Collections.ArrayList
    foreach ($rtAm in $rtAttrMatches) {
        $rtK = $rtAm.Groups[1].Value
        $rtV = $rtAm.Groups[2].Value
        $rtEntry = [PSCustomObject]@{ Key = $rtK; Value = $rtV }
        if ($rtK -eq "xmlns" -or $rtK.StartsWith("xmlns:")) {
            $rtNsAttrs.Add($rtEntry) | Out-Null
        } else {
            $rtNormAttrs.Add($rtEntry) | Out-Null
        }
    }

    $rtNsSorted = $rtNsAttrs | Sort-Object { Get-LocalName $_.Key }
    $rtNormSorted = $rtNormAttrs | Sort-Object { Get-LocalName $_.Key }

    $rtSb = New-Object System.Text.StringBuilder
    $rtSb.Append("<") | Out-Null
    $rtSb.Append($rtTagName) | Out-Null
    foreach ($rtE in $rtNsSorted) {
        $rtSb.Append(' ') | Out-Null
        $rtSb.Append($rtE.Key) | Out-Null
        $rtSb.Append('="') | Out-Null
        $rtSb.Append($rtE.Value) | Out-Null
        $rtSb.Append('"') | Out-Null
    }
    foreach ($rtE in $rtNormSorted) {
        $rtSb.Append(' ') | Out-Null
        $rtSb.Append($rtE.Key) | Out-Null
        $rtSb.Append('="') | Out-Null
        $rtSb.Append($rtE.Value) | Out-Null
        $rtSb.Append('"') | Out-Null
    }
    $rtSb.Append($rtSelfClose) | Out-Null
    $rtSb.Append(">") | Out-Null
    return $rtSb.ToString()
}

function Reorder-XmlAttributes {
    param($rxXmlText)

    $rxTagPattern = '<([\w:.\-]+)((?:\s+[^<>]*?)?)(/?)>'
    $rxMatches = [regex]::Matches($rxXmlText, $rxTagPattern)
    $rxSb = New-Object System.Text.StringBuilder
    $rxPos = 0
    foreach ($rxM in $rxMatches) {
        $rxTagName = $rxM.Groups[1].Value
        if ($rxTagName.StartsWith("?") -or $rxTagName.StartsWith("!")) {
            continue
        }
        $rxChunk = $rxXmlText.Substring($rxPos, $rxM.Index - $rxPos)
        $rxSb.Append($rxChunk) | Out-Null
        $rxNewTag = Reorder-Tag $rxM
        $rxSb.Append($rxNewTag) | Out-Null
        $rxPos = $rxM.Index + $rxM.Length
    }
    $rxSb.Append($rxXmlText.Substring($rxPos)) | Out-Null
    return $rxSb.ToString()
}

function Reorder-Part {
    param($rpFullText, $rpPartName)

    $rpNameMarker = 'pkg:name="' + $rpPartName + '"'
    $rpNameIdx = $rpFullText.IndexOf($rpNameMarker)
    if ($rpNameIdx -lt 0) {
        return $rpFullText
    }

    $rpDataStartMarker = "<pkg:xmlData>"
    $rpDataEndMarker = "</pkg:xmlData>"
    $rpDataMarkerIdx = $rpFullText.IndexOf($rpDataStartMarker, $rpNameIdx)
    $rpInnerStart = $rpDataMarkerIdx + $rpDataStartMarker.Length
    $rpInnerEnd = $rpFullText.IndexOf($rpDataEndMarker, $rpInnerStart)

    $rpInner = $rpFullText.Substring($rpInnerStart, $rpInnerEnd - $rpInnerStart)
    $rpNewInner = Reorder-XmlAttributes $rpInner

    $rpBefore = $rpFullText.Substring(0, $rpInnerStart)
    $rpAfter = $rpFullText.Substring($rpInnerEnd)

    return $rpBefore + $rpNewInner + $rpAfter
}

$d = $word.ActiveDocument
$xml = $d.WordOpenXML
$xml = Reorder-Part $xml "/word/document.xml"
$xml = Reorder-Part $xml "/word/styles.xml"
$d.WordOpenXML = $xml
Write-Output "done"
